# Auto-generated cell updates for Masamune_Profits workbook (ALC..WVR sheets)
# Updates currentAveragePrice / Leve price / profit columns (H:N) per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 36653
$ws.Range("I98").Value = 1409.4706
$ws.Range("J98").Value = 156481
$ws.Range("K98").Value = 1409.4706
$ws.Range("L98").Value = 156481
$ws.Range("M98").Value = 88.5293999999999
$ws.Range("N98").Value = -159477
$ws.Range("H122").Value = 36653
$ws.Range("I122").Value = 1409.4706
$ws.Range("J122").Value = 156481
$ws.Range("K122").Value = 4228.4118
$ws.Range("L122").Value = 469443
$ws.Range("M122").Value = -1778.4118
$ws.Range("N122").Value = -474343
$ws.Range("H126").Value = 46768
$ws.Range("J126").Value = 46768
$ws.Range("L126").Value = 46768
$ws.Range("N126").Value = -56648
$ws.Range("H137").Value = 5378.2256
$ws.Range("I137").Value = 1299.5385
$ws.Range("K137").Value = 3898.6155
$ws.Range("M137").Value = -1348.6155
$ws.Range("H138").Value = 1516.0515
$ws.Range("J138").Value = 2100
$ws.Range("L138").Value = 6300
$ws.Range("N138").Value = -16580

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10802.583
$ws.Range("I32").Value = 9484.804
$ws.Range("J32").Value = 29251.5
$ws.Range("K32").Value = 9484.804
$ws.Range("L32").Value = 29251.5
$ws.Range("M32").Value = -9197.804
$ws.Range("N32").Value = -29825.5
$ws.Range("H61").Value = 1634.5
$ws.Range("I61").Value = 1227.9333
$ws.Range("J61").Value = 2312.111
$ws.Range("K61").Value = 1227.9333
$ws.Range("L61").Value = 2312.111
$ws.Range("M61").Value = -1015.9333
$ws.Range("N61").Value = -2736.111
$ws.Range("H123").Value = 43810.25
$ws.Range("J123").Value = 43810.25
$ws.Range("L123").Value = 43810.25
$ws.Range("N123").Value = -53610.25
$ws.Range("H125").Value = 46195.8
$ws.Range("J125").Value = 46195.8
$ws.Range("L125").Value = 46195.8
$ws.Range("N125").Value = -56035.8
$ws.Range("H132").Value = 8334771.5
$ws.Range("I132").Value = 12196106
$ws.Range("J132").Value = 2419.6843
$ws.Range("K132").Value = 36588318
$ws.Range("L132").Value = 7259.0529
$ws.Range("M132").Value = -36585788
$ws.Range("N132").Value = -12319.0529
$ws.Range("H136").Value = 1634.5
$ws.Range("I136").Value = 1227.9333
$ws.Range("J136").Value = 2312.111
$ws.Range("K136").Value = 3683.7999
$ws.Range("L136").Value = 6936.333
$ws.Range("M136").Value = -1133.7999
$ws.Range("N136").Value = -12036.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 38281
$ws.Range("J92").Value = 38281
$ws.Range("L92").Value = 38281
$ws.Range("N92").Value = -43273
$ws.Range("H106").Value = 47984
$ws.Range("J106").Value = 47984
$ws.Range("L106").Value = 47984
$ws.Range("N106").Value = -50508
$ws.Range("H122").Value = 40641.6
$ws.Range("J122").Value = 40641.6
$ws.Range("L122").Value = 40641.6
$ws.Range("N122").Value = -50441.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2944.2556
$ws.Range("I31").Value = 957
$ws.Range("K31").Value = 957
$ws.Range("M31").Value = -662
$ws.Range("H34").Value = 2944.2556
$ws.Range("I34").Value = 957
$ws.Range("K34").Value = 957
$ws.Range("M34").Value = -755
$ws.Range("H82").Value = 40590.5
$ws.Range("J82").Value = 40590.5
$ws.Range("L82").Value = 40590.5
$ws.Range("N82").Value = -41312.5
$ws.Range("H85").Value = 40590.5
$ws.Range("J85").Value = 40590.5
$ws.Range("L85").Value = 40590.5
$ws.Range("N85").Value = -43086.5
$ws.Range("H92").Value = 44461.332
$ws.Range("J92").Value = 44461.332
$ws.Range("L92").Value = 44461.332
$ws.Range("N92").Value = -49453.332
$ws.Range("H100").Value = 40913
$ws.Range("J100").Value = 40913
$ws.Range("L100").Value = 40913
$ws.Range("N100").Value = -43077
$ws.Range("H106").Value = 30819.9
$ws.Range("J106").Value = 32622.375
$ws.Range("L106").Value = 32622.375
$ws.Range("N106").Value = -35146.375
$ws.Range("H122").Value = 60890.25
$ws.Range("I122").Value = 71470.88
$ws.Range("J122").Value = 933.3333
$ws.Range("K122").Value = 214412.64
$ws.Range("L122").Value = 2799.9999
$ws.Range("M122").Value = -211962.64
$ws.Range("N122").Value = -7699.9999
$ws.Range("H124").Value = 34263.6
$ws.Range("J124").Value = 34263.6
$ws.Range("L124").Value = 34263.6
$ws.Range("N124").Value = -39173.6
$ws.Range("H125").Value = 44996
$ws.Range("J125").Value = 44996
$ws.Range("L125").Value = 44996
$ws.Range("N125").Value = -49916
$ws.Range("H131").Value = 41896
$ws.Range("J131").Value = 41896
$ws.Range("L131").Value = 41896
$ws.Range("N131").Value = -51976

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1178.2
$ws.Range("I4").Value = 96.333336
$ws.Range("J4").Value = 2801
$ws.Range("K4").Value = 289.000008
$ws.Range("L4").Value = 8403
$ws.Range("M4").Value = -177.000008
$ws.Range("N4").Value = -8627
$ws.Range("H92").Value = 800
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1152
$ws.Range("N92").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1409.6923
$ws.Range("I113").Value = 1385.7142
$ws.Range("K113").Value = 1385.7142
$ws.Range("M113").Value = 784.2858000000001
$ws.Range("H118").Value = 39306
$ws.Range("J118").Value = 39306
$ws.Range("L118").Value = 39306
$ws.Range("N118").Value = -42620
$ws.Range("H120").Value = 39317
$ws.Range("J120").Value = 39317
$ws.Range("L120").Value = 39317
$ws.Range("N120").Value = -48993
$ws.Range("H125").Value = 44318
$ws.Range("J125").Value = 44318
$ws.Range("L125").Value = 44318
$ws.Range("N125").Value = -49238
$ws.Range("H127").Value = 46254.668
$ws.Range("J127").Value = 46254.668
$ws.Range("L127").Value = 46254.668
$ws.Range("N127").Value = -56174.668
$ws.Range("H130").Value = 53984
$ws.Range("J130").Value = 53984
$ws.Range("L130").Value = 53984
$ws.Range("N130").Value = -64024

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2768.4285
$ws.Range("I7").Value = 2133.25
$ws.Range("J7").Value = 4801
$ws.Range("K7").Value = 2133.25
$ws.Range("L7").Value = 4801
$ws.Range("M7").Value = -2021.25
$ws.Range("N7").Value = -5025
$ws.Range("H109").Value = 35152
$ws.Range("J109").Value = 35152
$ws.Range("L109").Value = 35152
$ws.Range("N109").Value = -37926
$ws.Range("H117").Value = 43384
$ws.Range("J117").Value = 43384
$ws.Range("L117").Value = 43384
$ws.Range("N117").Value = -52562
$ws.Range("H123").Value = 39421
$ws.Range("J123").Value = 39421
$ws.Range("L123").Value = 39421
$ws.Range("N123").Value = -49221
$ws.Range("H126").Value = 2768.4285
$ws.Range("I126").Value = 2133.25
$ws.Range("J126").Value = 4801
$ws.Range("K126").Value = 6399.75
$ws.Range("L126").Value = 14403
$ws.Range("M126").Value = -3929.75
$ws.Range("N126").Value = -19343
$ws.Range("H129").Value = 45411
$ws.Range("J129").Value = 45411
$ws.Range("L129").Value = 45411
$ws.Range("N129").Value = -55411
$ws.Range("H131").Value = 33254
$ws.Range("J131").Value = 33254
$ws.Range("L131").Value = 33254
$ws.Range("N131").Value = -43334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 31708.5
$ws.Range("I27").Value = 21000
$ws.Range("J27").Value = 42417
$ws.Range("K27").Value = 21000
$ws.Range("L27").Value = 42417
$ws.Range("M27").Value = -20931
$ws.Range("N27").Value = -42555
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H97").Value = 39572
$ws.Range("J97").Value = 39572
$ws.Range("L97").Value = 39572
$ws.Range("N97").Value = -41554
$ws.Range("H102").Value = 42337
$ws.Range("J102").Value = 42337
$ws.Range("L102").Value = 42337
$ws.Range("N102").Value = -48827
$ws.Range("H109").Value = 39373
$ws.Range("J109").Value = 39373
$ws.Range("L109").Value = 39373
$ws.Range("N109").Value = -42147
$ws.Range("H113").Value = 533.97144
$ws.Range("I113").Value = 518.2273
$ws.Range("J113").Value = 560.61536
$ws.Range("K113").Value = 1554.6819
$ws.Range("L113").Value = 1681.84608
$ws.Range("M113").Value = 615.3181
$ws.Range("N113").Value = -6021.84608
$ws.Range("H115").Value = 37456.332
$ws.Range("J115").Value = 37456.332
$ws.Range("L115").Value = 37456.332
$ws.Range("N115").Value = -40590.332
$ws.Range("H126").Value = 1029.5385
$ws.Range("I126").Value = 1073.6666
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 3220.9998
$ws.Range("L126").Value = 1500
$ws.Range("M126").Value = -750.9998000000001
$ws.Range("N126").Value = -6440
$ws.Range("H128").Value = 48218.75
$ws.Range("J128").Value = 48218.75
$ws.Range("L128").Value = 48218.75
$ws.Range("N128").Value = -58178.75
$ws.Range("H129").Value = 42429
$ws.Range("J129").Value = 42429
$ws.Range("L129").Value = 42429
$ws.Range("N129").Value = -52429
$ws.Range("H136").Value = 239952.23
$ws.Range("I136").Value = 297880.6
$ws.Range("J136").Value = 1802.3334
$ws.Range("K136").Value = 893641.7999999999
$ws.Range("L136").Value = 5407.0002
$ws.Range("M136").Value = -891091.7999999999
$ws.Range("N136").Value = -10507.0002
